$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row  = 99
$prev = $row - 1

# Plain numeric cells (B:F) -- set directly.
$ws.Cells.Item($row, 2).Value = 0   # volume
$ws.Cells.Item($row, 3).Value = 1   # high
$ws.Cells.Item($row, 4).Value = 1   # low
$ws.Cells.Item($row, 5).Value = 1   # open
$ws.Cells.Item($row, 6).Value = 1   # close

# G (adj_close) and H (ticker) are text cells stored as shared strings in
# this sheet ("1" and "YKY.MI", identical to the row above). Copying the
# row above's G:H cells (values + types) guarantees they land as shared
# text strings rather than being reinterpreted as numbers.
$ws.Range("G$prev`:H$prev").Copy()
$ws.Range("G$row`:H$row").PasteSpecial(-4104)   # xlPasteAll
$ws.Application.CutCopyMode = $false

# A (date) needs the same date number-format/style as the column above;
# copy just that format over, then set this row's own date value.
$ws.Range("A$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)          # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 45450.2916666667
